$wb = $excel.ActiveWorkbook

# Sheet ALC, row 86
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 2814.7222
$ws.Range("I86").Value = 2783.3333
$ws.Range("J86").Value = 2877.5
$ws.Range("K86").Value = 2783.3333
$ws.Range("L86").Value = 2877.5
$ws.Range("M86").Value = -1660.3333
$ws.Range("N86").Value = -5123.5

# Sheet ALC, row 89
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 2814.7222
$ws.Range("I89").Value = 2783.3333
$ws.Range("J89").Value = 2877.5
$ws.Range("K89").Value = 13916.6665
$ws.Range("L89").Value = 14387.5
$ws.Range("M89").Value = -8300.666499999999
$ws.Range("N89").Value = -25619.5

# Sheet ALC, row 92
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 8465.214
$ws.Range("I92").Value = 8465.214
$ws.Range("K92").Value = 8465.214
$ws.Range("M92").Value = -7217.214

# Sheet ALC, row 111
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 2390
$ws.Range("J111").Value = 2339.8
$ws.Range("L111").Value = 7019.400000000001
$ws.Range("N111").Value = -13153.4

# Sheet ALC, row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 13336326
$ws.Range("I116").Value = 25002286
$ws.Range("K116").Value = 25002286
$ws.Range("M116").Value = -24998844

# Sheet ALC, row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1284.7931
$ws.Range("I137").Value = 1178.7142
$ws.Range("J137").Value = 1563.25
$ws.Range("K137").Value = 3536.1426
$ws.Range("L137").Value = 4689.75
$ws.Range("M137").Value = -986.1425999999997
$ws.Range("N137").Value = -9789.75

# Sheet ALC, row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2713.3215
$ws.Range("I138").Value = 1586.25
$ws.Range("J138").Value = 4742.05
$ws.Range("K138").Value = 4758.75
$ws.Range("L138").Value = 14226.15
$ws.Range("M138").Value = 381.25
$ws.Range("N138").Value = -24506.15

# Sheet ARM, row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10347.076
$ws.Range("I32").Value = 10410.136
$ws.Range("K32").Value = 10410.136
$ws.Range("M32").Value = -10123.136

# Sheet ARM, row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2763.9
$ws.Range("I61").Value = 1997.5454
$ws.Range("K61").Value = 1997.5454
$ws.Range("M61").Value = -1785.5454

# Sheet ARM, row 117
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H117").Value = 25000
$ws.Range("J117").Value = 25000
$ws.Range("L117").Value = 25000
$ws.Range("N117").Value = -34178

# Sheet ARM, row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3880.4924
$ws.Range("I132").Value = 4882.343
$ws.Range("J132").Value = 2784.7188
$ws.Range("K132").Value = 14647.029
$ws.Range("L132").Value = 8354.1564
$ws.Range("M132").Value = -12117.029
$ws.Range("N132").Value = -13414.1564

# Sheet ARM, row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2763.9
$ws.Range("I136").Value = 1997.5454
$ws.Range("K136").Value = 5992.6362
$ws.Range("M136").Value = -3442.6362

# Sheet BSM, row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1680.4
$ws.Range("I99").Value = 901
$ws.Range("J99").Value = 2200
$ws.Range("K99").Value = 901
$ws.Range("L99").Value = 2200
$ws.Range("M99").Value = 597
$ws.Range("N99").Value = -5196

# Sheet BSM, row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2000.5902
$ws.Range("I134").Value = 1679.9791
$ws.Range("J134").Value = 3184.3845
$ws.Range("K134").Value = 5039.9373
$ws.Range("L134").Value = 9553.1535
$ws.Range("M134").Value = -2504.9373
$ws.Range("N134").Value = -14623.1535

# Sheet CRP, row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 789151.9
$ws.Range("I58").Value = 926903.6
$ws.Range("J58").Value = 1999.1428
$ws.Range("K58").Value = 926903.6
$ws.Range("L58").Value = 1999.1428
$ws.Range("M58").Value = -926700.6
$ws.Range("N58").Value = -2405.1428

# Sheet CRP, row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 789151.9
$ws.Range("I136").Value = 926903.6
$ws.Range("J136").Value = 1999.1428
$ws.Range("K136").Value = 2780710.8
$ws.Range("L136").Value = 5997.428400000001
$ws.Range("M136").Value = -2778160.8
$ws.Range("N136").Value = -11097.4284

# Sheet CUL, row 9
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 2000
$ws.Range("J9").Value = 2000
$ws.Range("L9").Value = 6000
$ws.Range("N9").Value = -6448

# Sheet CUL, row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 533.3125
$ws.Range("I107").Value = 544.3077
$ws.Range("J107").Value = 485.66666
$ws.Range("K107").Value = 1632.9231
$ws.Range("L107").Value = 1456.99998
$ws.Range("M107").Value = 287.0769
$ws.Range("N107").Value = -5296.999980000001

# Sheet CUL, row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1225.5714
$ws.Range("J122").Value = 1397.25
$ws.Range("L122").Value = 12575.25
$ws.Range("N122").Value = -17475.25

# Sheet CUL, row 125
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 4152.222
$ws.Range("I125").Value = 2030
$ws.Range("J125").Value = 4417.5
$ws.Range("K125").Value = 6090
$ws.Range("L125").Value = 13252.5
$ws.Range("M125").Value = -1170
$ws.Range("N125").Value = -23092.5

# Sheet CUL, row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 4836.6113
$ws.Range("I131").Value = 12891.125
$ws.Range("J131").Value = 2535.3215
$ws.Range("K131").Value = 38673.375
$ws.Range("L131").Value = 7605.9645
$ws.Range("M131").Value = -33633.375
$ws.Range("N131").Value = -17685.9645

# Sheet GSM, row 93
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 27150
$ws.Range("J93").Value = 27150
$ws.Range("L93").Value = 27150
$ws.Range("N93").Value = -30894

# Sheet GSM, row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3186.18
$ws.Range("I132").Value = 2997.907
$ws.Range("J132").Value = 4342.7144
$ws.Range("K132").Value = 8993.721000000001
$ws.Range("L132").Value = 13028.1432
$ws.Range("M132").Value = -6463.721000000001
$ws.Range("N132").Value = -18088.1432

# Sheet LTW, row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4749.6665
$ws.Range("I40").Value = 4699.6
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 4699.6
$ws.Range("L40").Value = 5000
$ws.Range("M40").Value = -4563.6
$ws.Range("N40").Value = -5272

# Sheet LTW, row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2975.5881
$ws.Range("I136").Value = 2240
$ws.Range("K136").Value = 6720
$ws.Range("M136").Value = -4170

# Sheet WVR, row 98
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H98").Value = 40589
$ws.Range("J98").Value = 40589
$ws.Range("L98").Value = 40589
$ws.Range("N98").Value = -46579

# Sheet WVR, row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 9767798
$ws.Range("I122").Value = 13890168
$ws.Range("J122").Value = 4467608
$ws.Range("K122").Value = 41670504
$ws.Range("L122").Value = 13402824
$ws.Range("M122").Value = -41668054
$ws.Range("N122").Value = -13407724

# Sheet WVR, row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1118.9429
$ws.Range("I136").Value = 1167.625
$ws.Range("J136").Value = 599.6667
$ws.Range("K136").Value = 3502.875
$ws.Range("L136").Value = 1799.0001
$ws.Range("M136").Value = -952.875
$ws.Range("N136").Value = -6899.0001
